$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AAP")

# Insert a new column before D. This shifts the existing D:K data right to
# E:L (matching the workbook's new dimension A5:L102 and the "spans"
# widening from 1:11 to 1:12 seen throughout the sheet).
$ws.Columns("D").Insert()

# The freshly inserted column D has no formatting; copy number formats
# (date style for the "Period Ending" rows, "#,##0" for data rows, bold
# font for the header row, etc.) from column E, which holds what used to
# be column D before the insert. Only the three blocks of rows that
# actually contained data cells in the old column D are touched (the
# blank separator rows 5/6/36/37/78/79 never had a D/E cell and must
# stay untouched).
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)

# Populate the new column D with the newly reported (fiscal 2018) figures.
$ws.Range("D7").Value = 43463
$ws.Range("D8").Value = 9580600
$ws.Range("D9").Value = 5354400
$ws.Range("D10").Value = 4226100
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 51500
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 8976300
$ws.Range("D18").Value = 604300
$ws.Range("D20").Value = 7600
$ws.Range("D21").Value = 850000
$ws.Range("D22").Value = 56600
$ws.Range("D23").Value = 555300
$ws.Range("D24").Value = 137100
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 418200
$ws.Range("D27").Value = 418200
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 5700
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -7600
$ws.Range("D33").Value = 423800
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 423800
$ws.Range("D38").Value = 43463
$ws.Range("D41").Value = 896500
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 625000
$ws.Range("D44").Value = 4362500
$ws.Range("D45").Value = 198400
$ws.Range("D46").Value = 6082500
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 1369000
$ws.Range("D49").Value = 1540800
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 48400
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 9040600
$ws.Range("D57").Value = 3172800
$ws.Range("D58").Value = "NA"
$ws.Range("D59").Value = 713200
$ws.Range("D60").Value = 3885900
$ws.Range("D61").Value = 1045700
$ws.Range("D62").Value = 558200
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 5489800
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 3326200
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 3550800
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43463
$ws.Range("D81").Value = 423800
$ws.Range("D83").Value = 238200
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 811000
$ws.Range("D91").Value = -193700
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -191800
$ws.Range("D96").Value = -17800
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -263900
$ws.Range("D101").Value = -5700
$ws.Range("D102").Value = 349600

Write-Output "done"
